$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T1").Value = "VQ_best"
$ws.Range("T2").Value = 1

$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("T5").Select()
